$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 14, shifting existing rows 14-17 down to 15-18
$ws.Rows.Item(14).Insert()

# Fill new row 14 with data (copy of old row 14 pattern, with updated values)
$ws.Cells.Item(14, 1).Value = 4
$ws.Cells.Item(14, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(14, 3).Value = "Los Lagos"
$ws.Cells.Item(14, 4).Value = 44873
$ws.Cells.Item(14, 5).Value = 10
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100101
$ws.Cells.Item(14, 8).Value = "Berries"
$ws.Cells.Item(14, 9).Value = 100101001
$ws.Cells.Item(14, 10).Value = "Arándano (blue)"
$ws.Cells.Item(14, 11).Value = "Sin especificar"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 400
$ws.Cells.Item(14, 14).Value = 7500
$ws.Cells.Item(14, 15).Value = 8000
$ws.Cells.Item(14, 16).Value = 7750
$ws.Cells.Item(14, 17).Value = "`$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(14, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(14, 19).Value = 5167
$ws.Cells.Item(14, 20).Value = 1.5
